$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "30.481.36"
$ws.Cells.Item(2, 5).Value = "  -0.79%  "

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.876.09"
$ws.Cells.Item(3, 5).Value = "  -0.74%  "

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.9996"
$ws.Cells.Item(4, 5).Value = "  -0.04%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "245.89"
$ws.Cells.Item(5, 5).Value = "  -1.05%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.9997"
$ws.Cells.Item(6, 5).Value = "  -0.11%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4711"
$ws.Cells.Item(7, 5).Value = "  -0.58%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -1.98%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.06516"
$ws.Cells.Item(9, 5).Value = "  -0.26%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "21.90"
$ws.Cells.Item(10, 5).Value = "  -0.60%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "101.25"
$ws.Cells.Item(11, 5).Value = "  +4.36%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.07812"
$ws.Cells.Item(12, 5).Value = "  -0.02%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "1.874.22"
$ws.Cells.Item(13, 5).Value = "  -0.86%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.7272"
$ws.Cells.Item(14, 5).Value = "  -1.09%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "5.179"
$ws.Cells.Item(15, 5).Value = "  -1.34%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "284.27"
$ws.Cells.Item(16, 5).Value = "  -0.14%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "30.464.90"
$ws.Cells.Item(17, 5).Value = "  -1.57%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "13.11"
$ws.Cells.Item(18, 5).Value = "  -0.81%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.9996"
$ws.Cells.Item(19, 5).Value = "  -0.11%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.000007492"
$ws.Cells.Item(20, 5).Value = "  -0.51%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "2.118.14"
$ws.Cells.Item(21, 5).Value = "  -0.99%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "5.336"
$ws.Cells.Item(22, 5).Value = "  +0.06%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.9992"
$ws.Cells.Item(23, 5).Value = "  -0.19%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "6.354"
$ws.Cells.Item(24, 5).Value = "  +1.55%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "9.060"
$ws.Cells.Item(25, 5).Value = "  -1.86%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "162.36"
$ws.Cells.Item(26, 5).Value = "  -1.18%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "19.00"
$ws.Cells.Item(27, 5).Value = "  +0.40%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "1.898"
$ws.Cells.Item(28, 5).Value = "  -1.48%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.09690"
$ws.Cells.Item(29, 5).Value = "  -0.49%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.323"
$ws.Cells.Item(30, 5).Value = "  -1.44%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.496"
$ws.Cells.Item(31, 5).Value = "  -0.15%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "4.230"
$ws.Cells.Item(32, 5).Value = "  -1.66%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "4.155"
$ws.Cells.Item(33, 5).Value = "  -0.83%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.04808"
$ws.Cells.Item(34, 5).Value = "  -1.15%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.127"
$ws.Cells.Item(35, 5).Value = "  +0.06%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.785"
$ws.Cells.Item(36, 5).Value = "  +2.23%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.6904"
$ws.Cells.Item(37, 5).Value = "  -1.12%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.01901"
$ws.Cells.Item(38, 5).Value = "  -0.29%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "2.849"
$ws.Cells.Item(39, 5).Value = "  +1.54%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "FraxShare"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "6.310"
$ws.Cells.Item(40, 5).Value = "  -1.39%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "Aave"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "75.81"
$ws.Cells.Item(41, 5).Value = "  -0.50%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.964"
$ws.Cells.Item(42, 5).Value = "  -2.56%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.4226"
$ws.Cells.Item(43, 5).Value = "  -0.83%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.9991"
$ws.Cells.Item(44, 5).Value = "  -0.19%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  -1.11%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "100.90"
$ws.Cells.Item(46, 5).Value = "  -0.67%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "9.732"
$ws.Cells.Item(47, 5).Value = "  +2.44%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "7.017"
$ws.Cells.Item(48, 5).Value = "  -0.32%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "35.05"
$ws.Cells.Item(49, 5).Value = "  -1.73%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.05759"
$ws.Cells.Item(50, 5).Value = "  +0.07%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "884.14"
$ws.Cells.Item(51, 5).Value = "  -3.85%  "
